$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G ("K") values regenerated: K used instead of Strike# (column G header is "K")
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 1
